# Rename the two logo pictures that appear (each twice - once per
# header/footer variant) in the document:
#   - the BTec logo picture:    name="image1.jpg" -> name="image2.jpg"
#   - the Pearson logo picture: name="image2.png" -> name="image1.png"
#
# Word's object model doesn't expose a "Name" property directly on
# InlineShape, so each inline picture is converted to a floating Shape
# (which does have a writable .Name), renamed, then converted back to
# an inline picture so the drawing stays <wp:inline>, matching the
# original layout.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-LogoPicture($range, $newName) {
    if ($range.InlineShapes.Count -gt 0) {
        $inlineShape = $range.InlineShapes.Item(1)
        $floatingShape = $inlineShape.ConvertToShape()
        $floatingShape.Name = $newName
        [void]$floatingShape.ConvertToInlineShape()
    }
}

# Headers: BTec_Logo-Orange picture, image1.jpg -> image2.jpg
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hdr = $sec.Headers.Item($i)
    Rename-LogoPicture $hdr.Range "image2.jpg"
}

# Footers: PearsonLogo picture, image2.png -> image1.png
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers.Item($i)
    Rename-LogoPicture $ftr.Range "image1.png"
}
